$wb = $excel.ActiveWorkbook

# RACP sheet: replace the formula in B2 with a plain (increased) value
$racp = $wb.Worksheets.Item("RACP")
$racp.Range("B2").Value = 160

# About sheet: move the active selection (cosmetic, matches saved cursor position)
$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("G17").Select()
